$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col3a1"
$ws.Range("C2").Value = "Mag"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.895758333333333
$ws.Range("H2").Value = 29.687275
$ws.Range("I2").Value = 0.009836335004010318
$ws.Range("J2").Value = 0.009836335004010316
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09045733333333333
$ws.Range("N2").Value = 0.271372
$ws.Range("O2").Value = 0.0574491187297735
$ws.Range("P2").Value = 0.0574491187297735
$ws.Range("Q2").Value = 0.8951439101444444
$ws.Range("R2").Value = 8.0562951913
$ws.Range("S2").Value = 0.0005650887775112158
$ws.Range("T2").Value = 0.0005650887775112157

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col3a1"
$ws.Range("C3").Value = "Mag"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.895758333333333
$ws.Range("H3").Value = 29.687275
$ws.Range("I3").Value = 0.009836335004010318
$ws.Range("J3").Value = 0.009836335004010316
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.484107
$ws.Range("N3").Value = 4.452321
$ws.Range("O3").Value = 0.9425508812702265
$ws.Range("P3").Value = 0.9425508812702265
$ws.Range("Q3").Value = 14.68636421280833
$ws.Range("R3").Value = 132.177277915275
$ws.Range("S3").Value = 0.009271246226499101
$ws.Range("T3").Value = 0.0092712462264991

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Col3a1"
$ws.Range("C4").Value = "Mag"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 950.6576336666667
$ws.Range("H4").Value = 2851.972901
$ws.Range("I4").Value = 0.9449490017724818
$ws.Range("J4").Value = 0.9449490017724816
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09045733333333333
$ws.Range("N4").Value = 0.271372
$ws.Range("O4").Value = 0.0574491187297735
$ws.Range("P4").Value = 0.0574491187297735
$ws.Range("Q4").Value = 85.99395445446356
$ws.Range("R4").Value = 773.945590090172
$ws.Range("S4").Value = 0.05428648739640825
$ws.Range("T4").Value = 0.05428648739640825

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col3a1"
$ws.Range("C5").Value = "Mag"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 950.6576336666667
$ws.Range("H5").Value = 2851.972901
$ws.Range("I5").Value = 0.9449490017724818
$ws.Range("J5").Value = 0.9449490017724816
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.484107
$ws.Range("N5").Value = 4.452321
$ws.Range("O5").Value = 0.9425508812702265
$ws.Range("P5").Value = 0.9425508812702265
$ws.Range("Q5").Value = 1410.877648728136
$ws.Range("R5").Value = 12697.89883855322
$ws.Range("S5").Value = 0.8906625143760735
$ws.Range("T5").Value = 0.8906625143760734

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Col3a1"
$ws.Range("C6").Value = "Mag"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 45.48781433333333
$ws.Range("H6").Value = 136.463443
$ws.Range("I6").Value = 0.04521466322350794
$ws.Range("J6").Value = 0.04521466322350793
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.09045733333333333
$ws.Range("N6").Value = 0.271372
$ws.Range("O6").Value = 0.0574491187297735
$ws.Range("P6").Value = 0.0574491187297735
$ws.Range("Q6").Value = 4.11470638375511
$ws.Range("R6").Value = 37.032357453796
$ws.Range("S6").Value = 0.002597542555854031
$ws.Range("T6").Value = 0.00259754255585403

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Col3a1"
$ws.Range("C7").Value = "Mag"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 45.48781433333333
$ws.Range("H7").Value = 136.463443
$ws.Range("I7").Value = 0.04521466322350794
$ws.Range("J7").Value = 0.04521466322350793
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.484107
$ws.Range("N7").Value = 4.452321
$ws.Range("O7").Value = 0.9425508812702265
$ws.Range("P7").Value = 0.9425508812702265
$ws.Range("Q7").Value = 67.50878366680033
$ws.Range("R7").Value = 607.579053001203
$ws.Range("S7").Value = 0.04261712066765391
$ws.Range("T7").Value = 0.0426171206676539

